# Insert two new price-report rows (889, 890) into the daily Fruta/Hortalizas
# "Limón" sheet for Vega Monumental Concepción, pushing all subsequent rows
# down by two (dimension grows from A1:T937 to A1:T939).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 889.
$ws.Rows.Item(889).Insert()
$ws.Rows.Item(889).Insert()

# New row 889: "1a amarillo"
$ws.Cells.Item(889, 1).Value = 11
$ws.Cells.Item(889, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(889, 3).Value = "Bíobío"
$ws.Cells.Item(889, 4).Value = "2023-12-07"
$ws.Cells.Item(889, 5).Value = 8
$ws.Cells.Item(889, 6).Value = "Fruta"
$ws.Cells.Item(889, 7).Value = 100102
$ws.Cells.Item(889, 8).Value = "Cítricos"
$ws.Cells.Item(889, 9).Value = 100102003
$ws.Cells.Item(889, 10).Value = "Limón"
$ws.Cells.Item(889, 11).Value = "Sin especificar"
$ws.Cells.Item(889, 12).Value = "1a amarillo"
$ws.Cells.Item(889, 13).Value = 180
$ws.Cells.Item(889, 14).Value = 12000
$ws.Cells.Item(889, 15).Value = 12000
$ws.Cells.Item(889, 16).Value = 12000
$ws.Cells.Item(889, 17).Value = "`$/malla 18 kilos"
$ws.Cells.Item(889, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(889, 19).Value = 667
$ws.Cells.Item(889, 20).Value = 18

# New row 890: "1a plateado"
$ws.Cells.Item(890, 1).Value = 11
$ws.Cells.Item(890, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(890, 3).Value = "Bíobío"
$ws.Cells.Item(890, 4).Value = "2023-12-07"
$ws.Cells.Item(890, 5).Value = 8
$ws.Cells.Item(890, 6).Value = "Fruta"
$ws.Cells.Item(890, 7).Value = 100102
$ws.Cells.Item(890, 8).Value = "Cítricos"
$ws.Cells.Item(890, 9).Value = 100102003
$ws.Cells.Item(890, 10).Value = "Limón"
$ws.Cells.Item(890, 11).Value = "Sin especificar"
$ws.Cells.Item(890, 12).Value = "1a plateado"
$ws.Cells.Item(890, 13).Value = 180
$ws.Cells.Item(890, 14).Value = 15000
$ws.Cells.Item(890, 15).Value = 15000
$ws.Cells.Item(890, 16).Value = 15000
$ws.Cells.Item(890, 17).Value = "`$/malla 18 kilos"
$ws.Cells.Item(890, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(890, 19).Value = 833
$ws.Cells.Item(890, 20).Value = 18
